$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1564
$ws1.Range("F5").Value = 166
$ws1.Range("F8").Value = 181
$ws1.Range("F9").Value = 759
$ws1.Range("F10").Value = 1055
$ws1.Range("F12").Value = 353
$ws1.Range("F13").Value = 64
$ws1.Range("F14").Value = 509
$ws1.Range("F15").Value = 22
$ws1.Range("F16").Value = 6518
$ws1.Range("F17").Value = 27
$ws1.Range("F18").Value = 86
$ws1.Range("F20").Value = 163
$ws1.Range("F22").Value = 15568
$ws1.Range("G22").Value = 60
$ws1.Range("F23").Value = 1539
$ws1.Range("F24").Value = 296
$ws1.Range("F25").Value = 151
$ws1.Range("F27").Value = 11113
$ws1.Range("F28").Value = 777
$ws1.Range("F29").Value = 4359
$ws1.Range("F30").Value = 253
$ws1.Range("F31").Value = 378
$ws1.Range("F32").Value = 23

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 346

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1564
$ws4.Range("F5").Value = 166
$ws4.Range("F7").Value = 346
$ws4.Range("F9").Value = 181
$ws4.Range("F10").Value = 759
$ws4.Range("F12").Value = 1055
$ws4.Range("F14").Value = 353
$ws4.Range("F15").Value = 64
$ws4.Range("F16").Value = 509
$ws4.Range("F18").Value = 22
$ws4.Range("F19").Value = 6519
$ws4.Range("F20").Value = 27
$ws4.Range("F21").Value = 86
$ws4.Range("F23").Value = 163
$ws4.Range("F26").Value = 15568
$ws4.Range("G26").Value = 60
$ws4.Range("F27").Value = 1539
$ws4.Range("F28").Value = 296
$ws4.Range("F29").Value = 151
$ws4.Range("F32").Value = 11113
$ws4.Range("F33").Value = 777
$ws4.Range("F34").Value = 4359
$ws4.Range("F35").Value = 253
$ws4.Range("F36").Value = 378
$ws4.Range("F37").Value = 23
